$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 64.8
$ws.Range("I11").Value = 64.8
$ws.Range("K11").Value = 64.8
$ws.Range("M11").Value = 75.2

$ws.Range("H17").Value = 2124923.8
$ws.Range("J17").Value = 2354080.5
$ws.Range("L17").Value = 7062241.5
$ws.Range("N17").Value = -7062577.5

$ws.Range("H19").Value = 1459.0476
$ws.Range("I19").Value = 2196.75
$ws.Range("K19").Value = 2196.75
$ws.Range("M19").Value = -2021.75

$ws.Range("H39").Value = 389
$ws.Range("I39").Value = 340.7857
$ws.Range("K39").Value = 1022.3571
$ws.Range("M39").Value = -726.3571000000001

$ws.Range("H43").Value = 8839
$ws.Range("I43").Value = 8811.6
$ws.Range("J43").Value = 8907.5
$ws.Range("K43").Value = 8811.6
$ws.Range("L43").Value = 8907.5
$ws.Range("M43").Value = -8742.6
$ws.Range("N43").Value = -9045.5

$ws.Range("H51").Value = 7130
$ws.Range("I51").Value = 8000
$ws.Range("J51").Value = 7084.2104
$ws.Range("K51").Value = 8000
$ws.Range("L51").Value = 7084.2104
$ws.Range("M51").Value = -7516
$ws.Range("N51").Value = -8052.2104

$ws.Range("H64").Value = 6881.3335
$ws.Range("J64").Value = 8497
$ws.Range("L64").Value = 8497
$ws.Range("N64").Value = -8993

$ws.Range("H67").Value = 6881.3335
$ws.Range("J67").Value = 8497
$ws.Range("L67").Value = 8497
$ws.Range("N67").Value = -10213

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H109").Value = 47485
$ws.Range("J109").Value = 47485
$ws.Range("L109").Value = 47485
$ws.Range("N109").Value = -50259

$ws.Range("H116").Value = 6844.4546
$ws.Range("I116").Value = 4617.857
$ws.Range("K116").Value = 4617.857
$ws.Range("M116").Value = -1175.857

$ws.Range("H137").Value = 8475.803
$ws.Range("I137").Value = 3727.6943
$ws.Range("K137").Value = 11183.0829
$ws.Range("M137").Value = -8633.082900000001

$ws.Range("H138").Value = 3755.4832
$ws.Range("J138").Value = 3834.946
$ws.Range("L138").Value = 11504.838
$ws.Range("N138").Value = -21784.838

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7223.1377
$ws.Range("I32").Value = 6981.5464
$ws.Range("K32").Value = 6981.5464
$ws.Range("M32").Value = -6694.5464

$ws.Range("H88").Value = 1849.65
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 1999.5883
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 1999.5883
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -2811.5883

$ws.Range("H91").Value = 1849.65
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 1999.5883
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 1999.5883
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -4807.588299999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7751.492
$ws.Range("I20").Value = 9117.232
$ws.Range("J20").Value = 4815.15
$ws.Range("K20").Value = 9117.232
$ws.Range("L20").Value = 4815.15
$ws.Range("M20").Value = -8870.232
$ws.Range("N20").Value = -5309.15

$ws.Range("H80").Value = 877.82355
$ws.Range("I80").Value = 621.5
$ws.Range("J80").Value = 1017.63635
$ws.Range("K80").Value = 621.5
$ws.Range("L80").Value = 1017.63635
$ws.Range("M80").Value = 376.5
$ws.Range("N80").Value = -3013.63635

$ws.Range("H83").Value = 877.82355
$ws.Range("I83").Value = 621.5
$ws.Range("J83").Value = 1017.63635
$ws.Range("K83").Value = 3107.5
$ws.Range("L83").Value = 5088.18175
$ws.Range("M83").Value = 1884.5
$ws.Range("N83").Value = -15072.18175

$ws.Range("H86").Value = 348328.38
$ws.Range("I86").Value = 910539.9399999999
$ws.Range("J86").Value = 4754.6665
$ws.Range("K86").Value = 910539.9399999999
$ws.Range("L86").Value = 4754.6665
$ws.Range("M86").Value = -909416.9399999999
$ws.Range("N86").Value = -7000.6665

$ws.Range("H89").Value = 348328.38
$ws.Range("I89").Value = 910539.9399999999
$ws.Range("J89").Value = 4754.6665
$ws.Range("K89").Value = 4552699.699999999
$ws.Range("L89").Value = 23773.3325
$ws.Range("M89").Value = -4547083.699999999
$ws.Range("N89").Value = -35005.3325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 486.8
$ws.Range("I94").Value = 312.33334
$ws.Range("J94").Value = 561.5714
$ws.Range("K94").Value = 312.33334
$ws.Range("L94").Value = 561.5714
$ws.Range("M94").Value = 138.66666
$ws.Range("N94").Value = -1463.5714

$ws.Range("H105").Value = 1300
$ws.Range("J105").Value = 1666.6666
$ws.Range("L105").Value = 1666.6666
$ws.Range("N105").Value = -5160.6666

$ws.Range("H134").Value = 2912.325
$ws.Range("I134").Value = 1439.08
$ws.Range("J134").Value = 5367.7334
$ws.Range("K134").Value = 4317.24
$ws.Range("L134").Value = 16103.2002
$ws.Range("M134").Value = -1782.24
$ws.Range("N134").Value = -21173.2002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 66.125
$ws.Range("I15").Value = 48.42857
$ws.Range("J15").Value = 190
$ws.Range("K15").Value = 145.28571
$ws.Range("L15").Value = 570
$ws.Range("M15").Value = -5.285709999999995
$ws.Range("N15").Value = -850

$ws.Range("H17").Value = 1258.9259
$ws.Range("J17").Value = 2025
$ws.Range("L17").Value = 6075
$ws.Range("N17").Value = -6413

$ws.Range("H32").Value = 1780.8
$ws.Range("J32").Value = 2300.3333
$ws.Range("L32").Value = 6900.999899999999
$ws.Range("N32").Value = -7466.999899999999

$ws.Range("H38").Value = 2905.6316
$ws.Range("J38").Value = 6277.75
$ws.Range("L38").Value = 18833.25
$ws.Range("N38").Value = -19527.25

$ws.Range("H80").Value = 42399.6
$ws.Range("I80").Value = 35666.332
$ws.Range("K80").Value = 106998.996
$ws.Range("M80").Value = -106062.996

$ws.Range("H83").Value = 42399.6
$ws.Range("I83").Value = 35666.332
$ws.Range("K83").Value = 320996.988
$ws.Range("M83").Value = -316316.988

$ws.Range("H113").Value = 2278.647
$ws.Range("J113").Value = 2217
$ws.Range("L113").Value = 6651
$ws.Range("N113").Value = -10991

$ws.Range("H122").Value = 7144120.5
$ws.Range("I122").Value = 1223.75
$ws.Range("J122").Value = 16667983
$ws.Range("K122").Value = 11013.75
$ws.Range("L122").Value = 150011847
$ws.Range("M122").Value = -8563.75
$ws.Range("N122").Value = -150016747

$ws.Range("H131").Value = 3663.8845
$ws.Range("J131").Value = 4488.7896
$ws.Range("L131").Value = 13466.3688
$ws.Range("N131").Value = -23546.3688

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 620
$ws.Range("I17").Value = 240
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 240
$ws.Range("L17").Value = 1000
$ws.Range("M17").Value = -72
$ws.Range("N17").Value = -1336

$ws.Range("H97").Value = 913.6667
$ws.Range("I97").Value = 758.5
$ws.Range("J97").Value = 1224
$ws.Range("K97").Value = 758.5
$ws.Range("L97").Value = 1224
$ws.Range("M97").Value = -262.5
$ws.Range("N97").Value = -2216

$ws.Range("H132").Value = 9446.429
$ws.Range("I132").Value = 8492.647000000001
$ws.Range("K132").Value = 25477.941
$ws.Range("M132").Value = -22947.941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2032.625
$ws.Range("I46").Value = 1118.3334
$ws.Range("K46").Value = 1118.3334
$ws.Range("M46").Value = -930.3334

$ws.Range("H82").Value = 2796.0715
$ws.Range("I82").Value = 2050.25
$ws.Range("J82").Value = 3094.4
$ws.Range("K82").Value = 2050.25
$ws.Range("L82").Value = 3094.4
$ws.Range("M82").Value = -1689.25
$ws.Range("N82").Value = -3816.4

$ws.Range("H85").Value = 2796.0715
$ws.Range("I85").Value = 2050.25
$ws.Range("J85").Value = 3094.4
$ws.Range("K85").Value = 2050.25
$ws.Range("L85").Value = 3094.4
$ws.Range("M85").Value = -802.25
$ws.Range("N85").Value = -5590.4

$ws.Range("H136").Value = 6071.9126
$ws.Range("I136").Value = 5632.4546
$ws.Range("K136").Value = 16897.3638
$ws.Range("M136").Value = -14347.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9957.166999999999
$ws.Range("I62").Value = 13230.667
$ws.Range("J62").Value = 8866
$ws.Range("K62").Value = 13230.667
$ws.Range("L62").Value = 8866
$ws.Range("M62").Value = -12606.667
$ws.Range("N62").Value = -10114

$ws.Range("H65").Value = 9957.166999999999
$ws.Range("I65").Value = 13230.667
$ws.Range("J65").Value = 8866
$ws.Range("K65").Value = 66153.33499999999
$ws.Range("L65").Value = 44330
$ws.Range("M65").Value = -63033.33499999999
$ws.Range("N65").Value = -50570

$ws.Range("H122").Value = 3702.8
$ws.Range("I122").Value = 3153.9092
$ws.Range("J122").Value = 5212.25
$ws.Range("K122").Value = 9461.7276
$ws.Range("L122").Value = 15636.75
$ws.Range("M122").Value = -7011.7276
$ws.Range("N122").Value = -20536.75

$ws.Range("H132").Value = 120862.125
$ws.Range("I132").Value = 171221.62
$ws.Range("K132").Value = 513664.86
$ws.Range("M132").Value = -511134.86
